$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-10-04"

# Update the header label in I1 to match the new "through" date
$ws.Range("I1").Value = "2022 (through 10-04)"

# Update October (row 11) value for the 2022 column
$ws.Range("I11").Value = 13

# Update Total (row 14) value for the 2022 column
$ws.Range("I14").Value = 1295
